$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the affected cells to keep Text formatting so the values stay as
# literal strings (matching the original inline-string cell type) instead of
# being auto-converted to numbers/percentages by Excel.
$cells = @("D2","E2","E3","D4","E4","D5","E5","D6","E6","E7","D8","E8","D9","E9","D10","E10","D11","E11","D12","E12","D13","E13","D14","E14","D15","E15","D16","E16","D17","E17","D18","E18","E19","D20","E20","D21","E21","D22","E22","D23","E23","D24","E24","D25","E25","D26","E26","D27","E27","D28","E28","D40","E40","D41","E41","E42","D43","E43","D44","E44","D45","E45","D47","E47")
foreach ($addr in $cells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values (Fri Jan 13 18:48:13 UTC 2023 symbol-list refresh).
$ws.Range("D2").Value = "287.56"
$ws.Range("E2").Value = "0.96%"
$ws.Range("E3").Value = "2.55%"
$ws.Range("D4").Value = "5.202"
$ws.Range("E4").Value = "2.13%"
$ws.Range("D5").Value = "0.06980"
$ws.Range("E5").Value = "4.89%"
$ws.Range("D6").Value = "7.444"
$ws.Range("E6").Value = "1.90%"
$ws.Range("E7").Value = "4.95%"
$ws.Range("D8").Value = "1.408"
$ws.Range("E8").Value = "4.00%"
$ws.Range("D9").Value = "0.9005"
$ws.Range("E9").Value = "-3.51%"
$ws.Range("D10").Value = "0.1609"
$ws.Range("E10").Value = "2.71%"
$ws.Range("D11").Value = "0.07521"
$ws.Range("E11").Value = "21.26%"
$ws.Range("D12").Value = "0.07721"
$ws.Range("E12").Value = "1.55%"
$ws.Range("D13").Value = "0.02944"
$ws.Range("E13").Value = "2.05%"
$ws.Range("D14").Value = "0.09016"
$ws.Range("E14").Value = "0.89%"
$ws.Range("D15").Value = "0.001572"
$ws.Range("E15").Value = "-0.79%"
$ws.Range("D16").Value = "0.0006516"
$ws.Range("E16").Value = "1.47%"
$ws.Range("D17").Value = "0.006070"
$ws.Range("E17").Value = "-0.81%"
$ws.Range("D18").Value = "3.473"
$ws.Range("E18").Value = "-0.26%"
$ws.Range("E19").Value = "0.18%"
$ws.Range("D20").Value = "0.3240"
$ws.Range("E20").Value = "1.32%"
$ws.Range("D21").Value = "0.1332"
$ws.Range("E21").Value = "2.27%"
$ws.Range("D22").Value = "4.000"
$ws.Range("E22").Value = "-1.19%"
$ws.Range("D23").Value = "0.1599"
$ws.Range("E23").Value = "5.19%"
$ws.Range("D24").Value = "0.04515"
$ws.Range("E24").Value = "1.24%"
$ws.Range("D25").Value = "0.001209"
$ws.Range("E25").Value = "2.65%"
$ws.Range("D26").Value = "0.004244"
$ws.Range("E26").Value = "-4.98%"
$ws.Range("D27").Value = "0.0001167"
$ws.Range("E27").Value = "-6.14%"
$ws.Range("D28").Value = "0.0001669"
$ws.Range("E28").Value = "3.67%"
$ws.Range("D40").Value = "0.04367"
$ws.Range("E40").Value = "4.78%"
$ws.Range("D41").Value = "0.006944"
$ws.Range("E41").Value = "3.67%"
$ws.Range("E42").Value = "0.18%"
$ws.Range("D43").Value = "0.002064"
$ws.Range("E43").Value = "2.75%"
$ws.Range("D44").Value = "0.01155"
$ws.Range("E44").Value = "0.67%"
$ws.Range("D45").Value = "0.00005831"
$ws.Range("E45").Value = "3.06%"
$ws.Range("D47").Value = "0.01307"
$ws.Range("E47").Value = "0.50%"
